# The sheet stores quarterly data in groups of 4 rows per year: A, B, C, D.
# This edit swaps the B-quarter and C-quarter rows (content, including the
# label in column A) within every 4-row year block, and removes the
# now-unused columns F ("家用电冰箱产销率") and G ("家用电冰箱销售量"),
# which were duplicates/derivations of columns B and E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$lastRow = 81

for ($r = $startRow; ($r + 3) -le $lastRow; $r += 4) {
    $rowB = $r + 1
    $rowC = $r + 2

    $rangeB = $ws.Range("A" + $rowB + ":E" + $rowB)
    $rangeC = $ws.Range("A" + $rowC + ":E" + $rowC)

    $valB = $rangeB.Value2()
    $valC = $rangeC.Value2()

    $rangeB.Value2 = $valC
    $rangeC.Value2 = $valB
}

# Remove columns F and G entirely (they duplicated B and E).
$ws.Range("F1:G81").Delete()
